$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/value updates (strings that are not ambiguous with numbers)
$ws.Range('D2').Value = '30.816.89'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '1.959.25'
$ws.Range('E3').Value = '  +3.68%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('E5').Value = '  +3.18%  '
$ws.Range('E6').Value = '  +37.94%  '
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('E8').Value = '  +11.98%  '
$ws.Range('E9').Value = '  +14.19%  '
$ws.Range('E10').Value = '  +5.64%  '
$ws.Range('E11').Value = '  +14.75%  '
$ws.Range('E12').Value = '  +2.87%  '
$ws.Range('E13').Value = '  +5.61%  '
$ws.Range('D14').Value = '1.944.10'
$ws.Range('E14').Value = '  +2.90%  '
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').Value = '30.824.67'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('E18').Value = '  +5.48%  '
$ws.Range('E19').Value = '  +2.80%  '
$ws.Range('E20').Value = '  +6.63%  '
$ws.Range('D21').Value = '2.195.16'
$ws.Range('E21').Value = '  +2.71%  '
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('E24').Value = '  +6.03%  '
$ws.Range('E25').Value = '  +4.78%  '
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('E27').Value = '  +2.91%  '
$ws.Range('E28').Value = '  +12.89%  '
$ws.Range('E29').Value = '  +26.94%  '
$ws.Range('E30').Value = '  +5.87%  '
$ws.Range('E31').Value = '  +1.24%  '
$ws.Range('E32').Value = '  +5.23%  '
$ws.Range('E33').Value = '  +5.44%  '
$ws.Range('E34').Value = '  +3.74%  '
$ws.Range('E35').Value = '  +6.43%  '
$ws.Range('E36').Value = '  +5.23%  '
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('E38').Value = '  +4.20%  '
$ws.Range('E39').Value = '  +3.73%  '
$ws.Range('E40').Value = '  +4.93%  '
$ws.Range('E41').Value = '  +2.46%  '
$ws.Range('E42').Value = '  +8.14%  '
$ws.Range('E43').Value = '  +1.87%  '
$ws.Range('E44').Value = '  +2.85%  '
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E47').Value = '  +3.82%  '
$ws.Range('E48').Value = '  +5.41%  '
$ws.Range('E49').Value = '  +2.67%  '
$ws.Range('E50').Value = '  +7.05%  '
$ws.Range('E51').Value = '  +1.44%  '

# Numeric-looking price strings: force Text storage, then restore the
# cell style to Normal so no visible formatting/style change remains
# (matches original file, which stores these as plain text cells).
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '0.9968'
$c.Style = "Normal"
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '251.70'
$c.Style = "Normal"
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.6485'
$c.Style = "Normal"
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.9965'
$c.Style = "Normal"
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3245'
$c.Style = "Normal"
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '25.36'
$c.Style = "Normal"
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.06866'
$c.Style = "Normal"
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.8335'
$c.Style = "Normal"
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.07975'
$c.Style = "Normal"
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '101.00'
$c.Style = "Normal"
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '5.375'
$c.Style = "Normal"
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '277.62'
$c.Style = "Normal"
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '13.76'
$c.Style = "Normal"
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.000007681'
$c.Style = "Normal"
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '5.630'
$c.Style = "Normal"
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.9992'
$c.Style = "Normal"
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '0.9968'
$c.Style = "Normal"
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '6.641'
$c.Style = "Normal"
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '9.505'
$c.Style = "Normal"
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '165.25'
$c.Style = "Normal"
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '19.48'
$c.Style = "Normal"
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.138'
$c.Style = "Normal"
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '0.1234'
$c.Style = "Normal"
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '1.555'
$c.Style = "Normal"
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '4.507'
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '4.373'
$c.Style = "Normal"
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.05038'
$c.Style = "Normal"
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.198'
$c.Style = "Normal"
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.7311'
$c.Style = "Normal"
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.710'
$c.Style = "Normal"
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.01972'
$c.Style = "Normal"
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '2.936'
$c.Style = "Normal"
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '6.513'
$c.Style = "Normal"
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '77.44'
$c.Style = "Normal"
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.4596'
$c.Style = "Normal"
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.030'
$c.Style = "Normal"
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.8466'
$c.Style = "Normal"
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.9974'
$c.Style = "Normal"
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '102.60'
$c.Style = "Normal"
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '9.991'
$c.Style = "Normal"
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '7.341'
$c.Style = "Normal"
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '35.99'
$c.Style = "Normal"
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.4210'
$c.Style = "Normal"
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '928.01'
$c.Style = "Normal"
